# Commit: "fully installed bot update buffer."
# Adds a new worksheet "Created update_history buffer" (a fresh day's
# timing-log entry, same layout/formulas as the preceding "Removed Index
# from Player " sheet) as the new last / active tab, and leaves the
# previously-active sheet's view in its post-navigation state.

$wb = $excel.ActiveWorkbook

# The template sheet carrying the day-to-day timing-log layout (headers,
# AVERAGE/SUM formulas, shared formulas, number formats) that every new
# day's entry is cloned from.
$template = $wb.Worksheets.Item($wb.Worksheets.Count)

# Clone it to the end of the tab strip -> becomes sheetId 10 / rId10,
# last position, and the active tab (matches activeTab 8 -> 9).
$template.Copy([System.Reflection.Missing]::Value, $template)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Created update_history buffer"

# New day's raw measurements (B/C/D = three timing runs; G2 = log date).
# E/B5:E6 stay formulas (AVERAGE / SUM / delta) copied verbatim from the
# template and recompute automatically from these inputs.
$newSheet.Range("B2").Value = 8.164
$newSheet.Range("C2").Value = 7.994
$newSheet.Range("D2").Value = 8.019

$newSheet.Range("B3").Value = 7.442
$newSheet.Range("C3").Value = 7.266
$newSheet.Range("D3").Value = 7.293

$newSheet.Range("B4").Value = 0.725
$newSheet.Range("C4").Value = 0.722
$newSheet.Range("D4").Value = 0.724

$newSheet.Range("G2").Value = 41802

# Restore the new sheet's intended selection state (matches the other
# freshly-active day sheets: activeCell D5).
$newSheet.Range("D5").Select()

# The previously-active sheet is no longer the selected tab; its window
# state moved on to a "select everything" view, and it now carries an
# explicit (portrait) page setup like its sibling sheets.
$template.Select()
$template.Range("A1:G6").Select()
$template.PageSetup().Orientation = 1

# Leave the new sheet as the active / visible tab.
$newSheet.Select()
